# Invert the sign of every numeric value in column E ("Block"),
# for all data rows (row 2 through the last used row), leaving
# blank/non-numeric cells untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)  # Column E
    $val = $cell.Value2
    if ($val -is [double] -or $val -is [int]) {
        $cell.Value2 = -1 * $val
    }
}
